# "Added Assign and resolve flag"
#
# 1. Rename the candidate name on SearchCases!A2 from "Meera Khan" to
#    "Janani Panchalingam" (shared string used elsewhere too).
# 2. Add a third worksheet "AssignFlag" after the existing two sheets,
#    with a small header/value table, and make it the active sheet/tab.
# 3. Restore/point the selection on each sheet to match the edited
#    workbook (SearchCases -> A2, AddNewCases -> C2, AssignFlag -> C2).

$wb = $excel.ActiveWorkbook

$wsSearch = $wb.Worksheets.Item("SearchCases")
$wsAddNew = $wb.Worksheets.Item("AddNewCases")

# --- 1. Update the candidate/searched name -------------------------------
$wsSearch.Range("A2").Value = "Janani Panchalingam"

# --- 2. Add the new "AssignFlag" worksheet at the end ---------------------
$wsFlag = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsFlag.Name = "AssignFlag"

$wsFlag.Range("A1").Value = "SearName"
$wsFlag.Range("B1").Value = "Reason"
$wsFlag.Range("C1").Value = "UnflagReason"

$wsFlag.Range("A2").Value = "Janani Panchalingam"
$wsFlag.Range("B2").Value = "Test Flag 1"
$wsFlag.Range("C2").Value = "UnflagReason"

# --- 3. Selections on each sheet ------------------------------------------
[void]$wsSearch.Range("A2").Select()
[void]$wsAddNew.Range("C2").Select()
[void]$wsFlag.Range("C2").Select()

# Make the new sheet the active tab (activeTab points to AssignFlag).
$wsFlag.Activate()
